$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45436
$ws.Range("D33").Value = 462
$ws.Range("D34").Value = 526
$ws.Range("D35").Value = 568
$ws.Range("D36").Value = 622
